$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17 (header pair, same look as row 15's sessionID/rating header) ---
$ws.Range("B15:C15").Copy()
$ws.Range("B17:C17").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B17").Value = "sessionID"
$ws.Range("C17").Value = "rating"

# --- Row 18 (data row, same look as row 16's data row) ---
$ws.Range("A16:C16").Copy()
$ws.Range("A18:C18").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A18").Value = 0
$ws.Range("B18").Value = 15
$ws.Range("C18").Value = 0

# --- Row 19 (header pair) ---
$ws.Range("B15:C15").Copy()
$ws.Range("B19:C19").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B19").Value = "sessionID"
$ws.Range("C19").Value = "rating"

# --- Row 20 (data row) ---
$ws.Range("A16:C16").Copy()
$ws.Range("A20:C20").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A20").Value = 0
$ws.Range("B20").Value = 17
$ws.Range("C20").Value = 5
